# Fruta / hortaliza, semanal
# Insert a new weekly record row before row 244 (pushing the existing
# rows 244-339 down to 245-340), then populate the newly inserted row
# with the new observation's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 244; Excel shifts rows 244:339 down
# to 245:340 and extends the sheet dimension to A1:R340.
$ws.Rows.Item(244).Insert()

# Fill the newly inserted row 244 with the new record.
$ws.Range("A244").Value = 5
$ws.Range("B244").Value = "Macroferia Regional de Talca"
$ws.Range("C244").Value = "Maule"
$ws.Range("D244").Value = 44784
$ws.Range("E244").Value = 7
$ws.Range("F244").Value = 100112003
$ws.Range("G244").Value = "Ajo"
$ws.Range("H244").Value = "Chino"
$ws.Range("I244").Value = "Primera"
$ws.Range("J244").Value = 300
$ws.Range("K244").Value = 28000
$ws.Range("L244").Value = 28000
$ws.Range("M244").Value = 28000
$ws.Range("N244").Value = "$/malla 10 kilos"
$ws.Range("O244").Value = "China"
$ws.Range("P244").Value = 2800
$ws.Range("Q244").Value = 10
$ws.Range("R244").Value = "Hortaliza"
